$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F45").Value = 66
$ws.Range("G45").Value = 12730.74
$ws.Range("F64").Value = 39
$ws.Range("G64").Value = 384.93
$ws.Range("B71").Value = 56279.78
$ws.Range("F100").Value = 36
$ws.Range("G100").Value = 6193.44
$ws.Range("F101").Value = 56
$ws.Range("G101").Value = 4599.84
$ws.Range("F105").Value = 222
$ws.Range("G105").Value = 18019.74
$ws.Range("F114").Value = 39
$ws.Range("G114").Value = 4489.68
$ws.Range("F115").Value = 125
$ws.Range("G115").Value = 7928.75
$ws.Range("F123").Value = 7
$ws.Range("G123").Value = 327.18
$ws.Range("F127").Value = 109
$ws.Range("G127").Value = 5537.2
$ws.Range("F141").Value = 399
$ws.Range("G141").Value = 7764.54
$ws.Range("B143").Value = 279683.53
$ws.Range("B213").Value = 53925
$ws.Range("B214").Value = 57756
$ws.Range("F234").Value = 27
$ws.Range("G234").Value = 2254.5
$ws.Range("F235").Value = 32
$ws.Range("G235").Value = 2701.44
$ws.Range("F237").Value = 10
$ws.Range("G237").Value = 1113.3
$ws.Range("B250").Value = 29034.45
$ws.Range("F256").Value = 56
$ws.Range("G256").Value = 1420.16
$ws.Range("B262").Value = 15761.98
$ws.Range("F272").Value = 2
$ws.Range("G272").Value = 168.38
$ws.Range("B273").Value = 168.38
$ws.Range("F296").Value = 76
$ws.Range("G296").Value = 6447.08
$ws.Range("B301").Value = 13803.57
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 65.73999999999999
$ws.Range("B319").Value = 31296.26
$ws.Range("F350").Value = 3
$ws.Range("G350").Value = 234.99
$ws.Range("F370").Value = 3
$ws.Range("G370").Value = 406.62
$ws.Range("F378").Value = 23
$ws.Range("G378").Value = 2327.6
$ws.Range("F381").Value = 4
$ws.Range("G381").Value = 1692.84
$ws.Range("B382").Value = 127224.16
$ws.Range("F387").Value = 246
$ws.Range("G387").Value = 25266.66
$ws.Range("F397").Value = 25
$ws.Range("G397").Value = 3016.25
$ws.Range("F405").Value = 110
$ws.Range("G405").Value = 15077.7
$ws.Range("F410").Value = 224
$ws.Range("G410").Value = 10671.36
$ws.Range("F420").Value = 289
$ws.Range("G420").Value = 15065.57
$ws.Range("F421").Value = 180
$ws.Range("G421").Value = 20152.8
$ws.Range("F422").Value = 174
$ws.Range("G422").Value = 22447.74
$ws.Range("F425").Value = 2
$ws.Range("G425").Value = 171.44
$ws.Range("F438").Value = 3
$ws.Range("G438").Value = 163.44
$ws.Range("F453").Value = 36
$ws.Range("G453").Value = 6695.64
$ws.Range("F456").Value = 25
$ws.Range("G456").Value = 3779.25
$ws.Range("B467").Value = 419743.12
$ws.Range("F492").Value = 9
$ws.Range("G492").Value = 324.9
$ws.Range("B500").Value = 1978.25
$ws.Range("F535").Value = 30
$ws.Range("G535").Value = 7545
$ws.Range("F537").Value = 69
$ws.Range("G537").Value = 2567.49
$ws.Range("B542").Value = 56551.96
$ws.Range("F556").Value = 60
$ws.Range("G556").Value = 2432.4
$ws.Range("F557").Value = 0
$ws.Range("G557").Value = 0
$ws.Range("B571").Value = 54050.12
$ws.Range("F587").Value = 656
$ws.Range("G587").Value = 8823.200000000001
$ws.Range("F591").Value = 231
$ws.Range("G591").Value = 4557.63
$ws.Range("F592").Value = 356
$ws.Range("G592").Value = 5849.08
$ws.Range("F593").Value = 333
$ws.Range("G593").Value = 4265.73
$ws.Range("F594").Value = 395
$ws.Range("G594").Value = 7793.35
$ws.Range("F595").Value = 456
$ws.Range("G595").Value = 3000.48
$ws.Range("F597").Value = 151
$ws.Range("G597").Value = 2938.46
$ws.Range("F603").Value = 313
$ws.Range("G603").Value = 5142.59
$ws.Range("B605").Value = 111486.32
$ws.Range("F607").Value = 64
$ws.Range("G607").Value = 2366.08
$ws.Range("B620").Value = 11551.19
$ws.Range("F670").Value = 439
$ws.Range("G670").Value = 3007.15
$ws.Range("F672").Value = 665
$ws.Range("G672").Value = 13200.25
$ws.Range("F674").Value = 327
$ws.Range("G674").Value = 5405.31
$ws.Range("B677").Value = 43284.27
$ws.Range("F697").Value = 36
$ws.Range("G697").Value = 7982.64
$ws.Range("B701").Value = 32744.44
$ws.Range("F724").Value = 34
$ws.Range("G724").Value = 751.74
$ws.Range("F725").Value = 130
$ws.Range("G725").Value = 5894.2
$ws.Range("F727").Value = 23
$ws.Range("G727").Value = 918.62
$ws.Range("B739").Value = 8556.92
$ws.Range("F748").Value = 5
$ws.Range("G748").Value = 605.85
$ws.Range("F757").Value = 31
$ws.Range("G757").Value = 2280.05
$ws.Range("F758").Value = 80
$ws.Range("G758").Value = 9826.4
$ws.Range("B761").Value = 55103.77
$ws.Range("F764").Value = 86
$ws.Range("G764").Value = 15309.72
$ws.Range("B770").Value = 69373.82000000001
$ws.Range("F796").Value = 288
$ws.Range("G796").Value = 4567.68
$ws.Range("F797").Value = 110
$ws.Range("G797").Value = 3642.1
$ws.Range("F798").Value = 193
$ws.Range("G798").Value = 8333.74
$ws.Range("B804").Value = 37509.77
$ws.Range("F834").Value = 26
$ws.Range("G834").Value = 2223
$ws.Range("F838").Value = 17
$ws.Range("G838").Value = 1366.46
$ws.Range("B839").Value = 4690.46
$ws.Range("F869").Value = 19
$ws.Range("G869").Value = 4479.06
$ws.Range("F875").Value = 9
$ws.Range("G875").Value = 1176.03
$ws.Range("F877").Value = 16
$ws.Range("G877").Value = 3927.52
$ws.Range("B880").Value = 28109.99
$ws.Range("F897").Value = 75
$ws.Range("G897").Value = 8312.25
$ws.Range("B901").Value = 52706.94
$ws.Range("F921").Value = 90
$ws.Range("G921").Value = 7073.1
$ws.Range("F923").Value = 88
$ws.Range("G923").Value = 9056.08
$ws.Range("B933").Value = 39292.96
$ws.Range("F935").Value = 2
$ws.Range("G935").Value = 215.16
$ws.Range("F936").Value = 71
$ws.Range("G936").Value = 2655.4
$ws.Range("F939").Value = 153
$ws.Range("G939").Value = 5722.2
$ws.Range("F941").Value = 131
$ws.Range("G941").Value = 4899.4
$ws.Range("B942").Value = 14377.74
$ws.Range("F974").Value = 6
$ws.Range("G974").Value = 5699.64
$ws.Range("B984").Value = 127488.25
$ws.Range("F999").Value = 1455
$ws.Range("G999").Value = 237325.05
$ws.Range("F1003").Value = 187
$ws.Range("G1003").Value = 12622.5
$ws.Range("B1005").Value = 279582.83
$ws.Range("B1012").Value = 2674588.49
$ws.Range("B1013").Value = 2674588.49
